# The deck currently carries two theme parts:
#   ppt/theme/theme1.xml -> "Integral"      (used by the slide master)
#   ppt/theme/theme2.xml -> "Office Theme"  (used by the notes master)
# The target edit swaps the two themes' colour schemes so that the slide
# master's theme becomes the stock "Office Theme" palette (and the notes
# master's theme becomes the "Integral" palette).
#
# The PowerPoint object model exposes theme colours through
# Theme.ThemeColorScheme.Colors(i).RGB (i = 1..12, in the fixed order
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) -- see e.g.
# $ppt.ActivePresentation.SlideMaster.Theme.ThemeColorScheme.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Target palette: the stock Office theme colours (RGB hex -> COM RGB int).
$colorScheme.Colors(1).RGB  = 0        # dk1      000000
$colorScheme.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$colorScheme.Colors(3).RGB  = 6968388  # dk2      44546A
$colorScheme.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$colorScheme.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$colorScheme.Colors(6).RGB  = 3243501  # accent2  ED7D31
$colorScheme.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$colorScheme.Colors(8).RGB  = 49407    # accent4  FFC000
$colorScheme.Colors(9).RGB  = 12874308 # accent5  4472C4
$colorScheme.Colors(10).RGB = 4697456  # accent6  70AD47
$colorScheme.Colors(11).RGB = 12673797 # hlink    0563C1
$colorScheme.Colors(12).RGB = 7491477  # folHlink 954F72
